# Fix Training Data Issue (#48)
# The "Date" column (BF) for each team row was stored using the sheet's
# original filename-style value (e.g. "6-14-2011-12"), which doesn't
# correctly represent the actual game date. Replace it with the correct
# ISO-style date string "2012-06-14" for every data row.
#
# NOTE: the target value looks like a date (YYYY-MM-DD), so if we just
# assign it to .Value, Excel's automatic type detection will convert it
# into a date serial number and slap a date NumberFormat on the cell.
# To keep it as plain text (matching the original cell's text storage)
# we force the cell's number format to Text ("@") before writing the
# value, then restore the cell style back to "Normal" so we don't leave
# a stray per-cell number format behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$correctDate = "2012-06-14"

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    $cell.NumberFormat = "@"
    $cell.Value = $correctDate
    $cell.Style = "Normal"
}
